$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1, matching the style of the other header cells
$ws.Range("F1").Value = "rejection-f"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 2 updates
$ws.Range("A2").Value = "even_MAG-GUT66378.fa"
$ws.Range("B2").Value = 0.06995035993880661
$ws.Range("C2").Value = 0.9300496400611934
$ws.Range("D2").Value = 0.9300496400611934
$ws.Range("F2").Value = "s__UMGS1370 sp900551135"

# Row 3 updates
$ws.Range("A3").Value = "even_MAG-GUT66382.fa"
$ws.Range("B3").Value = 0.04358882487501536
$ws.Range("C3").Value = 0.9564111751249846
$ws.Range("D3").Value = 0.9564111751249846
$ws.Range("F3").Value = "s__UMGS1370 sp900551135"

$wb.Save()
